$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change Log")

# Update the last existing "Notes" entry (row 7, column D) to drop the blank
# line before "Bugs" (Excel de-dupes/mutates the shared string in place).
$ws.Range("D7").Value = "Notes`n- Haven't tested what you added to alu_control or alu yet, but it all compiles so far`nBugs`n- "

# Fill in row 8 with the new change-log entry.
$ws.Range("A8").Value = "Night of 8/20/25"
$ws.Range("B8").Value = "Changes`n- MODIFIED: alu_control.vhd, alu.vhd`n- ADDED: Set on less than signed and unsigned                                                                                                                                                                                                                                        "
$ws.Range("D8").Value = "Notes`n- Compiles`nBugs`n- "

# Match the saved view state (scrolled/selected one row further down).
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("D9").Select()
